# Auto-generated script to update Leve profit/price columns (H-N) across all 8 sheets
# reflecting refreshed market-board data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 559.3077
$ws.Cells.Item(19, 9).Value = 545
$ws.Cells.Item(19, 10).Value = 565.6667
$ws.Cells.Item(19, 11).Value = 545
$ws.Cells.Item(19, 12).Value = 565.6667
$ws.Cells.Item(19, 13).Value = -370
$ws.Cells.Item(19, 14).Value = -915.6667
$ws.Cells.Item(28, 8).Value = 31621.559
$ws.Cells.Item(28, 9).Value = 40747.4
$ws.Cells.Item(28, 10).Value = 6272
$ws.Cells.Item(28, 11).Value = 40747.4
$ws.Cells.Item(28, 12).Value = 6272
$ws.Cells.Item(28, 13).Value = -40262.4
$ws.Cells.Item(28, 14).Value = -7242
$ws.Cells.Item(99, 8).Value = 523.1667
$ws.Cells.Item(99, 10).Value = 622.5
$ws.Cells.Item(99, 12).Value = 1867.5
$ws.Cells.Item(99, 14).Value = -4863.5
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 13).Value = $null
$ws.Cells.Item(113, 8).Value = 90913064
$ws.Cells.Item(113, 9).Value = 200002340
$ws.Cells.Item(113, 10).Value = 5334.3335
$ws.Cells.Item(113, 11).Value = 200002340
$ws.Cells.Item(113, 12).Value = 5334.3335
$ws.Cells.Item(113, 13).Value = -199999086
$ws.Cells.Item(113, 14).Value = -11842.3335
$ws.Cells.Item(118, 8).Value = 201
$ws.Cells.Item(118, 9).Value = 201
$ws.Cells.Item(118, 11).Value = 603
$ws.Cells.Item(118, 13).Value = 1054
$ws.Cells.Item(132, 8).Value = 5815.657
$ws.Cells.Item(132, 9).Value = 6031.815
$ws.Cells.Item(132, 11).Value = 18095.445
$ws.Cells.Item(132, 13).Value = -15565.445
$ws.Cells.Item(135, 8).Value = 333335500
$ws.Cells.Item(135, 9).Value = 1000000000
$ws.Cells.Item(135, 10).Value = 3249.5
$ws.Cells.Item(135, 11).Value = 9000000000
$ws.Cells.Item(135, 12).Value = 29245.5
$ws.Cells.Item(135, 13).Value = -8999997465
$ws.Cells.Item(135, 14).Value = -34315.5
$ws.Cells.Item(138, 8).Value = 2345.7576
$ws.Cells.Item(138, 9).Value = 1823.6923
$ws.Cells.Item(138, 11).Value = 5471.0769
$ws.Cells.Item(138, 13).Value = -331.0769
$ws.Cells.Item(141, 8).Value = 1801.5
$ws.Cells.Item(141, 9).Value = 1088.6666
$ws.Cells.Item(141, 10).Value = 2656.9
$ws.Cells.Item(141, 11).Value = 3265.9998
$ws.Cells.Item(141, 12).Value = 7970.700000000001
$ws.Cells.Item(141, 13).Value = 1914.0002
$ws.Cells.Item(141, 14).Value = -18330.7

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 1761.8889
$ws.Cells.Item(2, 9).Value = 1607.125
$ws.Cells.Item(2, 11).Value = 1607.125
$ws.Cells.Item(2, 13).Value = -1494.125
$ws.Cells.Item(32, 8).Value = 140540.88
$ws.Cells.Item(32, 9).Value = 142492.16
$ws.Cells.Item(32, 10).Value = 2000
$ws.Cells.Item(32, 11).Value = 142492.16
$ws.Cells.Item(32, 12).Value = 2000
$ws.Cells.Item(32, 13).Value = -142205.16
$ws.Cells.Item(32, 14).Value = -2574
$ws.Cells.Item(45, 8).Value = 2269
$ws.Cells.Item(45, 9).Value = 2031.1
$ws.Cells.Item(45, 11).Value = 2031.1
$ws.Cells.Item(45, 13).Value = -1654.1
$ws.Cells.Item(61, 8).Value = 2118.6365
$ws.Cells.Item(61, 9).Value = 2118.6365
$ws.Cells.Item(61, 11).Value = 2118.6365
$ws.Cells.Item(61, 13).Value = -1906.6365
$ws.Cells.Item(74, 8).Value = 3124.1943
$ws.Cells.Item(74, 9).Value = 3338.4285
$ws.Cells.Item(74, 11).Value = 3338.4285
$ws.Cells.Item(74, 13).Value = -2464.4285
$ws.Cells.Item(77, 8).Value = 3124.1943
$ws.Cells.Item(77, 9).Value = 3338.4285
$ws.Cells.Item(77, 11).Value = 16692.1425
$ws.Cells.Item(77, 13).Value = -12324.1425
$ws.Cells.Item(116, 8).Value = 1761.8889
$ws.Cells.Item(116, 9).Value = 1607.125
$ws.Cells.Item(116, 11).Value = 1607.125
$ws.Cells.Item(116, 13).Value = 686.875
$ws.Cells.Item(132, 8).Value = 47620844
$ws.Cells.Item(132, 9).Value = 90910720
$ws.Cells.Item(132, 10).Value = 1981.4
$ws.Cells.Item(132, 11).Value = 272732160
$ws.Cells.Item(132, 12).Value = 5944.200000000001
$ws.Cells.Item(132, 13).Value = -272729630
$ws.Cells.Item(132, 14).Value = -11004.2
$ws.Cells.Item(136, 8).Value = 2118.6365
$ws.Cells.Item(136, 9).Value = 2118.6365
$ws.Cells.Item(136, 11).Value = 6355.9095
$ws.Cells.Item(136, 13).Value = -3805.9095

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 1761.8889
$ws.Cells.Item(3, 9).Value = 1607.125
$ws.Cells.Item(3, 11).Value = 1607.125
$ws.Cells.Item(3, 13).Value = -1493.125
$ws.Cells.Item(20, 8).Value = 3922.3333
$ws.Cells.Item(20, 9).Value = 3193.875
$ws.Cells.Item(20, 10).Value = 9750
$ws.Cells.Item(20, 11).Value = 3193.875
$ws.Cells.Item(20, 12).Value = 9750
$ws.Cells.Item(20, 13).Value = -2946.875
$ws.Cells.Item(20, 14).Value = -10244
$ws.Cells.Item(74, 8).Value = 34124.375
$ws.Cells.Item(74, 10).Value = 31856.428
$ws.Cells.Item(74, 12).Value = 31856.428
$ws.Cells.Item(74, 14).Value = -33728.428
$ws.Cells.Item(77, 8).Value = 34124.375
$ws.Cells.Item(77, 10).Value = 31856.428
$ws.Cells.Item(77, 12).Value = 95569.284
$ws.Cells.Item(77, 14).Value = -104929.284
$ws.Cells.Item(86, 8).Value = 1877.0714
$ws.Cells.Item(86, 9).Value = 1735.4375
$ws.Cells.Item(86, 11).Value = 1735.4375
$ws.Cells.Item(86, 13).Value = -612.4375
$ws.Cells.Item(89, 8).Value = 1877.0714
$ws.Cells.Item(89, 9).Value = 1735.4375
$ws.Cells.Item(89, 11).Value = 8677.1875
$ws.Cells.Item(89, 13).Value = -3061.1875
$ws.Cells.Item(96, 8).Value = 14404
$ws.Cells.Item(96, 9).Value = 14404
$ws.Cells.Item(96, 11).Value = 14404
$ws.Cells.Item(96, 13).Value = -11658
$ws.Cells.Item(99, 8).Value = 1953.1765
$ws.Cells.Item(99, 9).Value = 1680.6
$ws.Cells.Item(99, 11).Value = 1680.6
$ws.Cells.Item(99, 13).Value = -182.5999999999999
$ws.Cells.Item(107, 8).Value = 12203969
$ws.Cells.Item(107, 9).Value = 4517.931
$ws.Cells.Item(107, 10).Value = 41685976
$ws.Cells.Item(107, 11).Value = 4517.931
$ws.Cells.Item(107, 12).Value = 41685976
$ws.Cells.Item(107, 13).Value = -2597.931
$ws.Cells.Item(107, 14).Value = -41689816
$ws.Cells.Item(134, 8).Value = 3080.5652
$ws.Cells.Item(134, 9).Value = 3097.65
$ws.Cells.Item(134, 10).Value = 2966.6667
$ws.Cells.Item(134, 11).Value = 9292.950000000001
$ws.Cells.Item(134, 12).Value = 8900.000100000001
$ws.Cells.Item(134, 13).Value = -6757.950000000001
$ws.Cells.Item(134, 14).Value = -13970.0001

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 3499.805
$ws.Cells.Item(31, 9).Value = 3109.5557
$ws.Cells.Item(31, 11).Value = 3109.5557
$ws.Cells.Item(31, 13).Value = -2814.5557
$ws.Cells.Item(34, 8).Value = 3499.805
$ws.Cells.Item(34, 9).Value = 3109.5557
$ws.Cells.Item(34, 11).Value = 3109.5557
$ws.Cells.Item(34, 13).Value = -2907.5557
$ws.Cells.Item(58, 8).Value = 2649.65
$ws.Cells.Item(58, 9).Value = 1788.8
$ws.Cells.Item(58, 11).Value = 1788.8
$ws.Cells.Item(58, 13).Value = -1585.8
$ws.Cells.Item(62, 8).Value = 8481.923000000001
$ws.Cells.Item(62, 9).Value = 8647
$ws.Cells.Item(62, 10).Value = 7931.6665
$ws.Cells.Item(62, 11).Value = 8647
$ws.Cells.Item(62, 12).Value = 7931.6665
$ws.Cells.Item(62, 13).Value = -8023
$ws.Cells.Item(62, 14).Value = -9179.666499999999
$ws.Cells.Item(65, 8).Value = 8481.923000000001
$ws.Cells.Item(65, 9).Value = 8647
$ws.Cells.Item(65, 10).Value = 7931.6665
$ws.Cells.Item(65, 11).Value = 43235
$ws.Cells.Item(65, 12).Value = 39658.3325
$ws.Cells.Item(65, 13).Value = -40115
$ws.Cells.Item(65, 14).Value = -45898.3325
$ws.Cells.Item(68, 8).Value = 56791.2
$ws.Cells.Item(68, 10).Value = 56791.2
$ws.Cells.Item(68, 12).Value = 56791.2
$ws.Cells.Item(68, 14).Value = -58289.2
$ws.Cells.Item(71, 8).Value = 56791.2
$ws.Cells.Item(71, 10).Value = 56791.2
$ws.Cells.Item(71, 12).Value = 170373.6
$ws.Cells.Item(71, 14).Value = -177861.6
$ws.Cells.Item(74, 8).Value = 36000
$ws.Cells.Item(74, 9).Value = 30500
$ws.Cells.Item(74, 10).Value = 43333.332
$ws.Cells.Item(74, 11).Value = 30500
$ws.Cells.Item(74, 12).Value = 43333.332
$ws.Cells.Item(74, 13).Value = -29626
$ws.Cells.Item(74, 14).Value = -45081.332
$ws.Cells.Item(77, 8).Value = 36000
$ws.Cells.Item(77, 9).Value = 30500
$ws.Cells.Item(77, 10).Value = 43333.332
$ws.Cells.Item(77, 11).Value = 91500
$ws.Cells.Item(77, 12).Value = 129999.996
$ws.Cells.Item(77, 13).Value = -87132
$ws.Cells.Item(77, 14).Value = -138735.996
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 14).Value = $null
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).Value = $null
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).Value = $null
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).Value = $null
$ws.Cells.Item(107, 8).Value = 665.92
$ws.Cells.Item(107, 9).Value = 716.15
$ws.Cells.Item(107, 11).Value = 716.15
$ws.Cells.Item(107, 13).Value = 1203.85
$ws.Cells.Item(132, 8).Value = 4635.8335
$ws.Cells.Item(132, 9).Value = 1672.3334
$ws.Cells.Item(132, 11).Value = 5017.0002
$ws.Cells.Item(132, 13).Value = -2487.0002
$ws.Cells.Item(135, 8).Value = 59999
$ws.Cells.Item(135, 10).Value = 59999
$ws.Cells.Item(135, 12).Value = 59999
$ws.Cells.Item(135, 14).Value = -70139
$ws.Cells.Item(136, 8).Value = 2649.65
$ws.Cells.Item(136, 9).Value = 1788.8
$ws.Cells.Item(136, 11).Value = 5366.4
$ws.Cells.Item(136, 13).Value = -2816.4

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(26, 8).Value = 387
$ws.Cells.Item(26, 10).Value = 599
$ws.Cells.Item(26, 12).Value = 1797
$ws.Cells.Item(26, 14).Value = -2373
$ws.Cells.Item(32, 8).Value = 1523.8
$ws.Cells.Item(32, 9).Value = 1934
$ws.Cells.Item(32, 10).Value = 1348
$ws.Cells.Item(32, 11).Value = 5802
$ws.Cells.Item(32, 12).Value = 4044
$ws.Cells.Item(32, 13).Value = -5519
$ws.Cells.Item(32, 14).Value = -4610
$ws.Cells.Item(124, 8).Value = 3500
$ws.Cells.Item(124, 9).Value = 0
$ws.Cells.Item(124, 10).Value = 3500
$ws.Cells.Item(124, 11).Value = 0
$ws.Cells.Item(124, 12).Value = 10500
$ws.Cells.Item(124, 13).Value = $null
$ws.Cells.Item(124, 14).Value = -20320
$ws.Cells.Item(132, 8).Value = 1419.4
$ws.Cells.Item(132, 10).Value = 1410.4445
$ws.Cells.Item(132, 12).Value = 12694.0005
$ws.Cells.Item(132, 14).Value = -17754.0005

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(14, 8).Value = 3751
$ws.Cells.Item(14, 9).Value = 3751
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 3751
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -3583
$ws.Cells.Item(14, 14).Value = $null
$ws.Cells.Item(26, 8).Value = 49999
$ws.Cells.Item(26, 10).Value = 49999
$ws.Cells.Item(26, 12).Value = 49999
$ws.Cells.Item(26, 14).Value = -50559
$ws.Cells.Item(50, 8).Value = 49999
$ws.Cells.Item(50, 10).Value = 49999
$ws.Cells.Item(50, 12).Value = 49999
$ws.Cells.Item(50, 14).Value = -50995
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).Value = $null

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 19233816
$ws.Cells.Item(7, 9).Value = 29414292
$ws.Cells.Item(7, 10).Value = 4027.111
$ws.Cells.Item(7, 11).Value = 29414292
$ws.Cells.Item(7, 12).Value = 4027.111
$ws.Cells.Item(7, 13).Value = -29414180
$ws.Cells.Item(7, 14).Value = -4251.111
$ws.Cells.Item(40, 8).Value = 3377.9048
$ws.Cells.Item(40, 9).Value = 2902.5625
$ws.Cells.Item(40, 10).Value = 4899
$ws.Cells.Item(40, 11).Value = 2902.5625
$ws.Cells.Item(40, 12).Value = 4899
$ws.Cells.Item(40, 13).Value = -2766.5625
$ws.Cells.Item(40, 14).Value = -5171
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 13).Value = $null
$ws.Cells.Item(126, 8).Value = 19233816
$ws.Cells.Item(126, 9).Value = 29414292
$ws.Cells.Item(126, 10).Value = 4027.111
$ws.Cells.Item(126, 11).Value = 88242876
$ws.Cells.Item(126, 12).Value = 12081.333
$ws.Cells.Item(126, 13).Value = -88240406
$ws.Cells.Item(126, 14).Value = -17021.333
$ws.Cells.Item(132, 8).Value = 6552.1904
$ws.Cells.Item(132, 9).Value = 3209.5
$ws.Cells.Item(132, 10).Value = 9591
$ws.Cells.Item(132, 11).Value = 9628.5
$ws.Cells.Item(132, 12).Value = 28773
$ws.Cells.Item(132, 13).Value = -7098.5
$ws.Cells.Item(132, 14).Value = -33833
$ws.Cells.Item(136, 8).Value = 6295.4346
$ws.Cells.Item(136, 9).Value = 4038.1904
$ws.Cells.Item(136, 10).Value = 29996.5
$ws.Cells.Item(136, 11).Value = 12114.5712
$ws.Cells.Item(136, 12).Value = 89989.5
$ws.Cells.Item(136, 13).Value = -9564.5712
$ws.Cells.Item(136, 14).Value = -95089.5

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 2646.5881
$ws.Cells.Item(122, 9).Value = 1894.8
$ws.Cells.Item(122, 11).Value = 5684.4
$ws.Cells.Item(122, 13).Value = -3234.4
$ws.Cells.Item(126, 8).Value = 1735.35
$ws.Cells.Item(126, 9).Value = 1655.9445
$ws.Cells.Item(126, 11).Value = 4967.833500000001
$ws.Cells.Item(126, 13).Value = -2497.833500000001
$ws.Cells.Item(132, 8).Value = 839068.2
$ws.Cells.Item(132, 9).Value = 1434459.4
$ws.Cells.Item(132, 11).Value = 4303378.199999999
$ws.Cells.Item(132, 13).Value = -4300848.199999999
$ws.Cells.Item(133, 8).Value = 62247.75
$ws.Cells.Item(133, 10).Value = 62247.75
$ws.Cells.Item(133, 12).Value = 62247.75
$ws.Cells.Item(133, 14).Value = -72367.75
